$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the forecast cells (C2, E2) entirely - the forecaster no longer
# produces a one-year-ahead forecast anchored on the first observation.
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: clear C3 (same bug fix as row 2) and update E3 with corrected value.
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 2.64925711235009

# Row 4: updated value (recalculated after bug fix)
$ws.Range("C4").Value = 2.533533936850585

# Row 5: updated value
$ws.Range("E5").Value = 2.332261646026246

# Row 6: updated values
$ws.Range("C6").Value = 1.21254482274098
$ws.Range("E6").Value = 1.839804681163337

# Row 7: updated value
$ws.Range("E7").Value = 0.6705904529405782

# Row 8: updated values
$ws.Range("C8").Value = 0.4712609263772816
$ws.Range("E8").Value = 0.8520644823059031

# Row 11: updated values
$ws.Range("C11").Value = 4.109890522944326
$ws.Range("E11").Value = 3.628019428949014

# Row 15: updated value
$ws.Range("E15").Value = 3.933586883651397

# Row 16: updated values
$ws.Range("C16").Value = 2.777797690741446
$ws.Range("E16").Value = 2.073300717643911

# Row 17: updated value
$ws.Range("E17").Value = 1.589741018019186

# Row 18: updated value
$ws.Range("C18").Value = -1.432689847121826

# Row 19: updated values
$ws.Range("C19").Value = 2.033479419175155
$ws.Range("E19").Value = 1.562315774899048
